$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H19").Value = 25005100
$wsALC.Range("I19").Value = 41668000
$wsALC.Range("J19").Value = 10750
$wsALC.Range("K19").Value = 41668000
$wsALC.Range("L19").Value = 10750
$wsALC.Range("M19").Value = -41667825
$wsALC.Range("N19").Value = -11100
$wsALC.Range("H129").Value = 279386.6
$wsALC.Range("I129").Value = 377.9565
$wsALC.Range("J129").Value = 773017.25
$wsALC.Range("K129").Value = 1133.8695
$wsALC.Range("L129").Value = 2319051.75
$wsALC.Range("M129").Value = 3866.1305
$wsALC.Range("N129").Value = -2329051.75
$wsALC.Range("H137").Value = 14547283
$wsALC.Range("I137").Value = 1115.2593
$wsALC.Range("J137").Value = 39093940
$wsALC.Range("K137").Value = 3345.7779
$wsALC.Range("L137").Value = 117281820
$wsALC.Range("M137").Value = -795.7779
$wsALC.Range("N137").Value = -117286920
$wsALC.Range("H140").Value = 68600
$wsALC.Range("J140").Value = 68600
$wsALC.Range("L140").Value = 68600
$wsALC.Range("N140").Value = -78960

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H32").Value = 5163080
$wsARM.Range("I32").Value = 7846.4175
$wsARM.Range("J32").Value = 27788826
$wsARM.Range("K32").Value = 7846.4175
$wsARM.Range("L32").Value = 27788826
$wsARM.Range("M32").Value = -7559.4175
$wsARM.Range("N32").Value = -27789400
$wsARM.Range("H61").Value = 2233940.2
$wsARM.Range("I61").Value = 1069214.6
$wsARM.Range("J61").Value = 9804658
$wsARM.Range("K61").Value = 1069214.6
$wsARM.Range("L61").Value = 9804658
$wsARM.Range("M61").Value = -1069002.6
$wsARM.Range("N61").Value = -9805082
$wsARM.Range("H132").Value = 8821379
$wsARM.Range("I132").Value = 10003121
$wsARM.Range("J132").Value = 4276219.5
$wsARM.Range("K132").Value = 30009363
$wsARM.Range("L132").Value = 12828658.5
$wsARM.Range("M132").Value = -30006833
$wsARM.Range("N132").Value = -12833718.5
$wsARM.Range("H135").Value = 32237
$wsARM.Range("I135").Value = 10000
$wsARM.Range("J135").Value = 34460.7
$wsARM.Range("K135").Value = 10000
$wsARM.Range("L135").Value = 34460.7
$wsARM.Range("M135").Value = -4930
$wsARM.Range("N135").Value = -44600.7
$wsARM.Range("H136").Value = 2233940.2
$wsARM.Range("I136").Value = 1069214.6
$wsARM.Range("J136").Value = 9804658
$wsARM.Range("K136").Value = 3207643.8
$wsARM.Range("L136").Value = 29413974
$wsARM.Range("M136").Value = -3205093.8
$wsARM.Range("N136").Value = -29419074
$wsARM.Range("H138").Value = 64476.332
$wsARM.Range("J138").Value = 64476.332
$wsARM.Range("L138").Value = 64476.332
$wsARM.Range("N138").Value = -74756.33199999999

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H86").Value = 1739.4736
$wsBSM.Range("I86").Value = 1788.5568
$wsBSM.Range("K86").Value = 1788.5568
$wsBSM.Range("M86").Value = -665.5568000000001
$wsBSM.Range("H89").Value = 1739.4736
$wsBSM.Range("I89").Value = 1788.5568
$wsBSM.Range("K89").Value = 8942.784
$wsBSM.Range("M89").Value = -3326.784

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H31").Value = 1044869.7
$wsCRP.Range("I31").Value = 1444.675
$wsCRP.Range("J31").Value = 3131719.8
$wsCRP.Range("K31").Value = 1444.675
$wsCRP.Range("L31").Value = 3131719.8
$wsCRP.Range("M31").Value = -1149.675
$wsCRP.Range("N31").Value = -3132309.8
$wsCRP.Range("H34").Value = 1044869.7
$wsCRP.Range("I34").Value = 1444.675
$wsCRP.Range("J34").Value = 3131719.8
$wsCRP.Range("K34").Value = 1444.675
$wsCRP.Range("L34").Value = 3131719.8
$wsCRP.Range("M34").Value = -1242.675
$wsCRP.Range("N34").Value = -3132123.8
$wsCRP.Range("H99").Value = 42400
$wsCRP.Range("I99").Value = 75000
$wsCRP.Range("J99").Value = 34250
$wsCRP.Range("K99").Value = 75000
$wsCRP.Range("L99").Value = 34250
$wsCRP.Range("M99").Value = -73502
$wsCRP.Range("N99").Value = -37246
$wsCRP.Range("H126").Value = 42400
$wsCRP.Range("I126").Value = 75000
$wsCRP.Range("J126").Value = 34250
$wsCRP.Range("K126").Value = 225000
$wsCRP.Range("L126").Value = 102750
$wsCRP.Range("M126").Value = -222530
$wsCRP.Range("N126").Value = -107690
$wsCRP.Range("H132").Value = 1533.9454
$wsCRP.Range("I132").Value = 1236.3334
$wsCRP.Range("J132").Value = 2873.2
$wsCRP.Range("K132").Value = 3709.0002
$wsCRP.Range("L132").Value = 8619.599999999999
$wsCRP.Range("M132").Value = -1179.0002
$wsCRP.Range("N132").Value = -13679.6
$wsCRP.Range("H134").Value = 873908.75
$wsCRP.Range("I134").Value = 4252.919
$wsCRP.Range("J134").Value = 4449160.5
$wsCRP.Range("K134").Value = 12758.757
$wsCRP.Range("L134").Value = 13347481.5
$wsCRP.Range("M134").Value = -10223.757
$wsCRP.Range("N134").Value = -13352551.5

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H69").Value = 47620844
$wsCUL.Range("I69").Value = 0
$wsCUL.Range("J69").Value = 47620844
$wsCUL.Range("K69").Value = 0
$wsCUL.Range("L69").Value = 142862532
$wsCUL.Range("M69").ClearContents()
$wsCUL.Range("N69").Value = -142864154
$wsCUL.Range("H72").Value = 47620844
$wsCUL.Range("I72").Value = 0
$wsCUL.Range("J72").Value = 47620844
$wsCUL.Range("K72").Value = 0
$wsCUL.Range("L72").Value = 428587596
$wsCUL.Range("M72").ClearContents()
$wsCUL.Range("N72").Value = -428595708
$wsCUL.Range("H107").Value = 801588.1
$wsCUL.Range("I107").Value = 2849199.2
$wsCUL.Range("J107").Value = 348.91306
$wsCUL.Range("K107").Value = 8547597.600000001
$wsCUL.Range("L107").Value = 1046.73918
$wsCUL.Range("M107").Value = -8545677.600000001
$wsCUL.Range("N107").Value = -4886.73918
$wsCUL.Range("H113").Value = 3019.449
$wsCUL.Range("I113").Value = 458.2857
$wsCUL.Range("J113").Value = 4940.3213
$wsCUL.Range("K113").Value = 1374.8571
$wsCUL.Range("L113").Value = 14820.9639
$wsCUL.Range("M113").Value = 795.1428999999998
$wsCUL.Range("N113").Value = -19160.9639
$wsCUL.Range("H120").Value = 333333340
$wsCUL.Range("I120").Value = 333333340
$wsCUL.Range("J120").Value = 0
$wsCUL.Range("K120").Value = 1000000020
$wsCUL.Range("L120").Value = 0
$wsCUL.Range("M120").Value = -999995182
$wsCUL.Range("N120").ClearContents()
$wsCUL.Range("H131").Value = 7057160.5
$wsCUL.Range("I131").Value = 31312866
$wsCUL.Range("J131").Value = 955.16364
$wsCUL.Range("K131").Value = 93938598
$wsCUL.Range("L131").Value = 2865.49092
$wsCUL.Range("M131").Value = -93933558
$wsCUL.Range("N131").Value = -12945.49092

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H132").Value = 16928376
$wsGSM.Range("I132").Value = 15477990
$wsGSM.Range("J132").Value = 22729922
$wsGSM.Range("K132").Value = 46433970
$wsGSM.Range("L132").Value = 68189766
$wsGSM.Range("M132").Value = -46431440
$wsGSM.Range("N132").Value = -68194826

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H132").Value = 2509561.2
$wsLTW.Range("I132").Value = 3666697.2
$wsLTW.Range("J132").Value = 2433.5
$wsLTW.Range("K132").Value = 11000091.6
$wsLTW.Range("L132").Value = 7300.5
$wsLTW.Range("M132").Value = -10997561.6
$wsLTW.Range("N132").Value = -12360.5
$wsLTW.Range("H136").Value = 2268840.8
$wsLTW.Range("I136").Value = 2316004
$wsLTW.Range("J136").Value = 5000
$wsLTW.Range("K136").Value = 6948012
$wsLTW.Range("L136").Value = 15000
$wsLTW.Range("M136").Value = -6945462
$wsLTW.Range("N136").Value = -20100

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H122").Value = 2925348.8
$wsWVR.Range("I122").Value = 3832607
$wsWVR.Range("J122").Value = 1961.1111
$wsWVR.Range("K122").Value = 11497821
$wsWVR.Range("L122").Value = 5883.3333
$wsWVR.Range("M122").Value = -11495371
$wsWVR.Range("N122").Value = -10783.3333
$wsWVR.Range("H132").Value = 995039.75
$wsWVR.Range("I132").Value = 3788.037
$wsWVR.Range("J132").Value = 2269506.2
$wsWVR.Range("K132").Value = 11364.111
$wsWVR.Range("L132").Value = 6808518.600000001
$wsWVR.Range("M132").Value = -8834.110999999999
$wsWVR.Range("N132").Value = -6813578.600000001
$wsWVR.Range("H140").Value = 57925
$wsWVR.Range("J140").Value = 57925
$wsWVR.Range("L140").Value = 57925
$wsWVR.Range("N140").Value = -68285
